$wb = $excel.ActiveWorkbook

# --- Mmusculs sheet: record the 19 Sept 2015 mouse archive, assembled up through the 06 Sept download ---
$wsMouse = $wb.Worksheets.Item("Mmusculs")
$wsMouse.Range("B6").Value = '$killdevil:/ms/home/s/b/sbiswas/transcriptome_compression/Mmusculus/NCBI_SRA_Mmusculus_successful_downloads_19Sept2015_download_compiled_23Sept2015.tgz'
$wsMouse.Range("A6").Select()

# --- Athaliana sheet: move the selection cursor (stays the active/selected tab) ---
$wsPlant = $wb.Worksheets.Item("Athaliana")
$wsPlant.Activate()
$wsPlant.Range("B6").Select()
